$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2943368.2
$ws.Range("I132").Value = 3280386.5
$ws.Range("K132").Value = 9841159.5
$ws.Range("M132").Value = -9838629.5
$ws.Range("H135").Value = 978.6539
$ws.Range("I135").Value = 565.5263
$ws.Range("J135").Value = 2100
$ws.Range("K135").Value = 5089.736699999999
$ws.Range("L135").Value = 18900
$ws.Range("M135").Value = -2554.736699999999
$ws.Range("N135").Value = -23970
$ws.Range("H137").Value = 3129913.8
$ws.Range("I137").Value = 4353069
$ws.Range("J137").Value = 4073.2222
$ws.Range("K137").Value = 13059207
$ws.Range("L137").Value = 12219.6666
$ws.Range("M137").Value = -13056657
$ws.Range("N137").Value = -17319.6666
$ws.Range("H138").Value = 6432.4185
$ws.Range("I138").Value = 2941.04
$ws.Range("J138").Value = 7628.0957
$ws.Range("K138").Value = 8823.119999999999
$ws.Range("L138").Value = 22884.2871
$ws.Range("M138").Value = -3683.119999999999
$ws.Range("N138").Value = -33164.2871

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18296.766
$ws.Range("I32").Value = 15431.786
$ws.Range("K32").Value = 15431.786
$ws.Range("M32").Value = -15144.786
$ws.Range("H37").Value = 8036
$ws.Range("J37").Value = 16038
$ws.Range("L37").Value = 16038
$ws.Range("N37").Value = -16584
$ws.Range("H45").Value = 1117.0698
$ws.Range("I45").Value = 917.90247
$ws.Range("J45").Value = 5200
$ws.Range("K45").Value = 917.90247
$ws.Range("L45").Value = 5200
$ws.Range("M45").Value = -540.90247
$ws.Range("N45").Value = -5954
$ws.Range("H61").Value = 2759.2368
$ws.Range("I61").Value = 1993.5217
$ws.Range("J61").Value = 3933.3333
$ws.Range("K61").Value = 1993.5217
$ws.Range("L61").Value = 3933.3333
$ws.Range("M61").Value = -1781.5217
$ws.Range("N61").Value = -4357.3333
$ws.Range("H64").Value = 29200
$ws.Range("J64").Value = 29200
$ws.Range("L64").Value = 29200
$ws.Range("N64").Value = -29696
$ws.Range("H67").Value = 29200
$ws.Range("J67").Value = 29200
$ws.Range("L67").Value = 29200
$ws.Range("N67").Value = -30916
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140
$ws.Range("H136").Value = 2759.2368
$ws.Range("I136").Value = 1993.5217
$ws.Range("J136").Value = 3933.3333
$ws.Range("K136").Value = 5980.5651
$ws.Range("L136").Value = 11799.9999
$ws.Range("M136").Value = -3430.5651
$ws.Range("N136").Value = -16899.9999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 48748.184
$ws.Range("I86").Value = 5266.6665
$ws.Range("J86").Value = 65053.75
$ws.Range("K86").Value = 5266.6665
$ws.Range("L86").Value = 65053.75
$ws.Range("M86").Value = -4143.6665
$ws.Range("N86").Value = -67299.75
$ws.Range("H89").Value = 48748.184
$ws.Range("I89").Value = 5266.6665
$ws.Range("J89").Value = 65053.75
$ws.Range("K89").Value = 26333.3325
$ws.Range("L89").Value = 325268.75
$ws.Range("M89").Value = -20717.3325
$ws.Range("N89").Value = -336500.75
$ws.Range("H99").Value = 4428.0586
$ws.Range("I99").Value = 3251.3076
$ws.Range("J99").Value = 8252.5
$ws.Range("K99").Value = 3251.3076
$ws.Range("L99").Value = 8252.5
$ws.Range("M99").Value = -1753.3076
$ws.Range("N99").Value = -11248.5
$ws.Range("H105").Value = 2291.9092
$ws.Range("I105").Value = 2093.8462
$ws.Range("J105").Value = 2578
$ws.Range("K105").Value = 2093.8462
$ws.Range("L105").Value = 2578
$ws.Range("M105").Value = -346.8462
$ws.Range("N105").Value = -6072
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 29337.75
$ws.Range("J137").Value = 28999
$ws.Range("L137").Value = 28999
$ws.Range("N137").Value = -39199
$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2044051.4
$ws.Range("I31").Value = 2633479.8
$ws.Range("J31").Value = 7844.909
$ws.Range("K31").Value = 2633479.8
$ws.Range("L31").Value = 7844.909
$ws.Range("M31").Value = -2633184.8
$ws.Range("N31").Value = -8434.909
$ws.Range("H34").Value = 2044051.4
$ws.Range("I34").Value = 2633479.8
$ws.Range("J34").Value = 7844.909
$ws.Range("K34").Value = 2633479.8
$ws.Range("L34").Value = 7844.909
$ws.Range("M34").Value = -2633277.8
$ws.Range("N34").Value = -8248.909
$ws.Range("H58").Value = 8623172
$ws.Range("I58").Value = 1232
$ws.Range("J58").Value = 27783040
$ws.Range("K58").Value = 1232
$ws.Range("L58").Value = 27783040
$ws.Range("M58").Value = -1029
$ws.Range("N58").Value = -27783446
$ws.Range("H60").Value = 9333.333000000001
$ws.Range("J60").Value = 9000
$ws.Range("L60").Value = 9000
$ws.Range("N60").Value = -10022
$ws.Range("H74").Value = 24648
$ws.Range("J74").Value = 24648
$ws.Range("L74").Value = 24648
$ws.Range("N74").Value = -26396
$ws.Range("H77").Value = 24648
$ws.Range("J77").Value = 24648
$ws.Range("L77").Value = 73944
$ws.Range("N77").Value = -82680
$ws.Range("H132").Value = 2234.3618
$ws.Range("I132").Value = 1514.3715
$ws.Range("J132").Value = 4334.3335
$ws.Range("K132").Value = 4543.1145
$ws.Range("L132").Value = 13003.0005
$ws.Range("M132").Value = -2013.1145
$ws.Range("N132").Value = -18063.0005
$ws.Range("H134").Value = 2135.62
$ws.Range("I134").Value = 1577.1714
$ws.Range("J134").Value = 3438.6667
$ws.Range("K134").Value = 4731.5142
$ws.Range("L134").Value = 10316.0001
$ws.Range("M134").Value = -2196.5142
$ws.Range("N134").Value = -15386.0001
$ws.Range("H136").Value = 8623172
$ws.Range("I136").Value = 1232
$ws.Range("J136").Value = 27783040
$ws.Range("K136").Value = 3696
$ws.Range("L136").Value = 83349120
$ws.Range("M136").Value = -1146
$ws.Range("N136").Value = -83354220
$ws.Range("H141").Value = 31008.334
$ws.Range("J141").Value = 32097.059
$ws.Range("L141").Value = 32097.059
$ws.Range("N141").Value = -42457.059

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1257.8
$ws.Range("I107").Value = 556
$ws.Range("J107").Value = 1635.6923
$ws.Range("K107").Value = 1668
$ws.Range("L107").Value = 4907.0769
$ws.Range("M107").Value = 252
$ws.Range("N107").Value = -8747.0769
$ws.Range("H131").Value = 1476.6833
$ws.Range("I131").Value = 5849.75
$ws.Range("J131").Value = 1164.3214
$ws.Range("K131").Value = 17549.25
$ws.Range("L131").Value = 3492.9642
$ws.Range("M131").Value = -12509.25
$ws.Range("N131").Value = -13572.9642

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 11500
$ws.Range("I57").Value = 11500
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 11500
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -10680
$ws.Range("N57").ClearContents()
$ws.Range("H97").Value = 1741.3334
$ws.Range("I97").Value = 1342.8572
$ws.Range("J97").Value = 2090
$ws.Range("K97").Value = 1342.8572
$ws.Range("L97").Value = 2090
$ws.Range("M97").Value = -846.8571999999999
$ws.Range("N97").Value = -3082
$ws.Range("H132").Value = 4735.5757
$ws.Range("I132").Value = 5118.857
$ws.Range("J132").Value = 4064.8333
$ws.Range("K132").Value = 15356.571
$ws.Range("L132").Value = 12194.4999
$ws.Range("M132").Value = -12826.571
$ws.Range("N132").Value = -17254.4999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2829.125
$ws.Range("I7").Value = 1324
$ws.Range("J7").Value = 3330.8333
$ws.Range("K7").Value = 1324
$ws.Range("L7").Value = 3330.8333
$ws.Range("M7").Value = -1212
$ws.Range("N7").Value = -3554.8333
$ws.Range("H61").Value = 71431960
$ws.Range("I61").Value = 111112180
$ws.Range("J61").Value = 7560
$ws.Range("K61").Value = 111112180
$ws.Range("L61").Value = 7560
$ws.Range("M61").Value = -111111978
$ws.Range("N61").Value = -7964
$ws.Range("H100").Value = 4004.5454
$ws.Range("I100").Value = 1750
$ws.Range("J100").Value = 4505.5557
$ws.Range("K100").Value = 1750
$ws.Range("L100").Value = 4505.5557
$ws.Range("M100").Value = -1209
$ws.Range("N100").Value = -5587.5557
$ws.Range("H113").Value = 71431960
$ws.Range("I113").Value = 111112180
$ws.Range("J113").Value = 7560
$ws.Range("K113").Value = 111112180
$ws.Range("L113").Value = 7560
$ws.Range("M113").Value = -111110010
$ws.Range("N113").Value = -11900
$ws.Range("H126").Value = 2829.125
$ws.Range("I126").Value = 1324
$ws.Range("J126").Value = 3330.8333
$ws.Range("K126").Value = 3972
$ws.Range("L126").Value = 9992.499899999999
$ws.Range("M126").Value = -1502
$ws.Range("N126").Value = -14932.4999
$ws.Range("H136").Value = 2946310.8
$ws.Range("I136").Value = 4170694.2
$ws.Range("J136").Value = 7790.5
$ws.Range("K136").Value = 12512082.6
$ws.Range("L136").Value = 23371.5
$ws.Range("M136").Value = -12509532.6
$ws.Range("N136").Value = -28471.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2708.6592
$ws.Range("I136").Value = 1972.5
$ws.Range("J136").Value = 4286.143
$ws.Range("K136").Value = 5917.5
$ws.Range("L136").Value = 12858.429
$ws.Range("M136").Value = -3367.5
$ws.Range("N136").Value = -17958.429
